$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 439.464458638125
$ws.Range("C2").Value = 514.43989045750004
$ws.Range("D2").Value = 437.90090132812503
$ws.Range("E2").Value = 522.57770713312505

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 432.99391268812496
$ws.Range("C3").Value = 525.7632940725
$ws.Range("D3").Value = 445.52344103999997
$ws.Range("E3").Value = 522.94383455249999

# Update the selection to match the new active range (B1:E3)
$ws.Range("B1:E3").Select()
